# Updated cryptos list on Wed Dec  6 14:56:14 UTC 2023 with GitHub Actions
#
# Refreshes the live snapshot in the "cryptos" worksheet: new Price (D)
# and Volume(1h) (E) readings for every coin, plus a rank swap between
# HuobiToken and NEARProtocol (rows 50-51) now that HuobiToken overtook
# NEARProtocol in the feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "44.091.04"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "2.255.48"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'229.65"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "'0.631"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "'63.18"
$ws.Range("E7").Value = "  +4.06%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.441"
$ws.Range("E9").Value = "  +9.87%  "
$ws.Range("E10").Value = "  +14.27%  "
$ws.Range("D11").Value = "'57.26"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "'25.80"
$ws.Range("E12").Value = "  +16.76%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "2.594.60"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "'15.58"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("E16").Value = "  +10.65%  "
$ws.Range("D17").Value = "'0.836"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").Value = "2.270.88"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "43.942.05"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("E20").Value = "  +8.88%  "
$ws.Range("D21").Value = "'73.00"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "'6.00"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "'250.60"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'9.97"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("D28").Value = "'3.22"
$ws.Range("E28").Value = "  +21.43%  "
$ws.Range("D29").Value = "'171.73"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "'20.69"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "'0.135"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("D32").Value = "'1.37"
$ws.Range("E32").Value = "  -6.12%  "
$ws.Range("D33").Value = "'0.123"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("D34").Value = "'0.0678"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").Value = "'4.70"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").Value = "'4.82"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Value = "'3.79"
$ws.Range("E37").Value = "  +6.81%  "
$ws.Range("D38").Value = "'6.61"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "'0.0257"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "'17.29"
$ws.Range("E42").Value = "  +8.36%  "
$ws.Range("D43").Value = "'8.21"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "'0.0961"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "'97.03"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'0.000209"
$ws.Range("E47").Value = "  -8.39%  "
$ws.Range("D48").Value = "'4.32"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "1.434.06"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.75"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.26"
$ws.Range("E51").Value = "  +2.84%  "
